$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Sheet1"

# New header row (columns A-K), replacing the old A-M layout
$headers = @("Description", "Remediation", "PowerShell Script", "Returned Value", "Default Value", "Expected Value", "Impact", "Likelihood", "Priority", "RiskRating", "References")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# New row 2 data values (columns A-K)
$row2 = @(
    "Enabling the default Teams DLP policy rule in Microsoft 365 helps protect an organization's sensitive information by preventing accidental sharing or leakage of that information in Teams conversations and channels.",
    "Use the PowerShell script to create a new DLPCompliancePolicy or review the policies existence and if they are enabled.",
    'New-DlpCompliancePolicy -Name "SSN Teams Policy" -Comment "SSN Teams Policy" -TeamsLocation All -Mode Enable',
    "Default Value:Enable",
    "Enable",
    "Enable",
    "2",
    "1",
    "Informational",
    "Low",
    "Learn about data loss prevention"
)
# Columns G ("Impact") and H ("Likelihood") hold numeric-looking text
# ("2" / "1"); force text format so they stay strings instead of numbers.
$textColumns = @(7, 8)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(2, $col)
    if ($textColumns -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $row2[$i]
}

# Remove the now-unused trailing columns (old L/M), shifting the sheet
# down to the new A1:K2 dimension
$ws.Columns("L:M").Delete()
